$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "21_01_2024"
$ws.Range("E2").Value = 2638
$ws.Range("E3").Value = 2032
$ws.Range("E4").Value = 3446
$ws.Range("E5").Value = 6784

$ws.Range("E6").Select() | Out-Null
